# Apply cryptos price/volume update (Sat Aug 31 14:51:26 UTC 2024 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.068.40"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.506.37"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'536.63"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'135.97"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "2.517.89"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("D12").Value = "'5.32"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "'0.347"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "2.940.67"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "'22.97"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "58.868.63"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "2.517.04"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'11.05"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "'4.25"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'322.94"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").Value = "'65.02"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'7.51"
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("D29").Value = "0.0₃0766"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'6.63"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.76"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'170.44"
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("E33").Value = "  +8.15%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'18.34"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'4.06"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "'1.53"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'36.88"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'284.23"
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.15"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'0.606"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "'129.85"
$ws.Range("E46").Value = "  +4.76%  "
$ws.Range("D47").Value = "'10.85"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'0.0922"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'0.0503"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'17.35"
$ws.Range("E51").Value = "  -0.07%  "
